$wb = $excel.ActiveWorkbook

$wsContracts  = $wb.Worksheets.Item("Contracts")
$wsActivities = $wb.Worksheets.Item("Activities")
$wsMedia      = $wb.Worksheets.Item("Media")

# Fill the new "Client" column (C) on the Contracts sheet for rows 2-28
$wsContracts.Range("C2:C28").Value = "Dubai Holding"

# Update each sheet's saved selection/active-cell state
$wsContracts.Activate()
$wsContracts.Range("C2:C28").Select()

$wsActivities.Activate()
$wsActivities.Range("E50").Select()

$wsMedia.Activate()
$wsMedia.Range("A2:A28").Select()

# Media ends up as the active/selected sheet (tab) once the file is saved
$wsMedia.Activate()
